$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.25067138671875
$ws.Range("B1").Value = 2.194236278533936
$ws.Range("C1").Value = 2.772485256195068
$ws.Range("D1").Value = 3.222027063369751
$ws.Range("E1").Value = 2.272112369537354
